$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.854.57"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.314.76"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0784"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").Value = "2.678.04"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").Value = "2.313.71"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").Value = "42.787.83"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("E19").Value = "  -6.19%  "

$ws.Range("E20").Value = "  +2.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.77%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("E28").Value = "  +14.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0698"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.92%  "

$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.26%  "

$ws.Range("D43").Value = "1.927.57"
$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "2.546.66"
$ws.Range("E48").Value = "  +1.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.29%  "
